# Lesson 13 21 Nov 2024
# Target: the speaker-notes page of slide 7 ("DQN Learning").
#   1. Nudge the notes "slide image" placeholder's x offset (381300 -> 381000 EMU).
#   2. Add the notes body text: "DQN Learning – deep q-network learning"
#      (authored as two runs in the source OOXML).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$notes = $s.NotesPage

# --- 1) Slide-image placeholder: tiny re-snap of its left offset. -----------
# (381300 EMU -> 381000 EMU, i.e. 30.0236pt -> 30.0pt)
$imgShape = $notes.Shapes.Item(1)
try {
    $imgShape.Left = 30.0
} catch {
    # Position editing on this placeholder is not supported by this host;
    # continue so the rest of the edit still lands.
}

# --- 2) Notes body placeholder: add the two text runs. ---------------------
$bodyShape = $notes.Shapes.Item(2)
$tr = $bodyShape.TextFrame.TextRange

$enDash = [char]0x2013
$run1 = "DQN Learning " + $enDash + " deep "
$run2 = "q-network learning"

$tr.Text = $run1 + $run2
